$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns keep their exact text representation (avoid numeric coercion)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.641.04"
$ws.Range("E2").Value = "  +2.79%  "
$ws.Range("D3").Value = "2.481.61"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "575.21"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").Value = "150.29"
$ws.Range("E6").Value = "  +5.24%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "0.541"
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("E9").Value = "  +4.97%  "
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "0.366"
$ws.Range("E11").Value = "  +4.64%  "
$ws.Range("D12").Value = "5.36"
$ws.Range("E12").Value = "  +3.03%  "
$ws.Range("D13").Value = "27.35"
$ws.Range("E13").Value = "  +6.31%  "
$ws.Range("D14").Value = "0.0000186"
$ws.Range("E14").Value = "  +7.65%  "
$ws.Range("D15").Value = "2.900.70"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").Value = "63.475.89"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").Value = "2.480.47"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").Value = "11.63"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("D19").Value = "7.28"
$ws.Range("E19").Value = "  +6.76%  "
$ws.Range("D20").Value = "4.27"
$ws.Range("E20").Value = "  +3.43%  "
$ws.Range("D21").Value = "329.77"
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").Value = "1.91"
$ws.Range("E23").Value = "  +9.89%  "
$ws.Range("D24").Value = "67.79"
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("D25").Value = "640.29"
$ws.Range("E25").Value = "  +14.74%  "
$ws.Range("E26").Value = "  +13.76%  "
$ws.Range("D27").Value = "8.88"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "1.54"
$ws.Range("E28").Value = "  +10.86%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.603.32"
$ws.Range("E29").Value = "  +2.79%  "
$ws.Range("D30").Value = "8.57"
$ws.Range("E30").Value = "  +4.92%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("D33").Value = "1.93"
$ws.Range("E33").Value = "  +3.70%  "
$ws.Range("D34").Value = "5.23"
$ws.Range("E34").Value = "  +10.04%  "
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  +4.36%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "0.388"
$ws.Range("E36").Value = "  +2.62%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "0.997"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "5.55"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("D39").Value = "18.98"
$ws.Range("E39").Value = "  +2.75%  "
$ws.Range("D40").Value = "1.86"
$ws.Range("E40").Value = "  +2.58%  "
$ws.Range("D41").Value = "147.73"
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("D42").Value = "2.65"
$ws.Range("E42").Value = "  +17.97%  "
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").Value = "152.50"
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("D45").Value = "3.80"
$ws.Range("E45").Value = "  +4.92%  "
$ws.Range("D46").Value = "0.0554"
$ws.Range("E46").Value = "  +5.48%  "
$ws.Range("D47").Value = "21.22"
$ws.Range("E47").Value = "  +7.42%  "
$ws.Range("D48").Value = "0.614"
$ws.Range("E48").Value = "  +3.69%  "
$ws.Range("D49").Value = "0.0240"
$ws.Range("E49").Value = "  +6.32%  "
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").Value = "0.745"
$ws.Range("E51").Value = "  +5.55%  "
